$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 173.82051
$ws.Range("I9").Value = 111.90625
$ws.Range("J9").Value = 456.85715
$ws.Range("K9").Value = 111.90625
$ws.Range("L9").Value = 456.85715
$ws.Range("M9").Value = 57.09375
$ws.Range("N9").Value = -794.85715

$ws.Range("H41").Value = 278.125
$ws.Range("J41").Value = 353.4
$ws.Range("L41").Value = 353.4
$ws.Range("N41").Value = -1233.4

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H98").Value = 1750.8182
$ws.Range("I98").Value = 1517.6666
$ws.Range("K98").Value = 1517.6666
$ws.Range("M98").Value = -19.66660000000002

$ws.Range("H122").Value = 1750.8182
$ws.Range("I122").Value = 1517.6666
$ws.Range("K122").Value = 4552.9998
$ws.Range("M122").Value = -2102.9998

$ws.Range("H132").Value = 8990
$ws.Range("I132").Value = 4247.952
$ws.Range("J132").Value = 21437.875
$ws.Range("K132").Value = 12743.856
$ws.Range("L132").Value = 64313.625
$ws.Range("M132").Value = -10213.856
$ws.Range("N132").Value = -69373.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15418.5
$ws.Range("I45").Value = 23755
$ws.Range("K45").Value = 23755
$ws.Range("M45").Value = -23378

$ws.Range("H46").Value = 8761.429
$ws.Range("J46").Value = 10028
$ws.Range("L46").Value = 10028
$ws.Range("N46").Value = -10666

$ws.Range("H61").Value = 8207.5
$ws.Range("I61").Value = 2332
$ws.Range("J61").Value = 14083
$ws.Range("K61").Value = 2332
$ws.Range("L61").Value = 14083
$ws.Range("M61").Value = -2120
$ws.Range("N61").Value = -14507

$ws.Range("H122").Value = 3111.5625
$ws.Range("I122").Value = 2127.5715
$ws.Range("K122").Value = 6382.7145
$ws.Range("M122").Value = -3932.7145

$ws.Range("H136").Value = 8207.5
$ws.Range("I136").Value = 2332
$ws.Range("J136").Value = 14083
$ws.Range("K136").Value = 6996
$ws.Range("L136").Value = 42249
$ws.Range("M136").Value = -4446
$ws.Range("N136").Value = -47349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3968.3928
$ws.Range("J20").Value = 3964.125
$ws.Range("L20").Value = 3964.125
$ws.Range("N20").Value = -4458.125

$ws.Range("H22").Value = 553.2
$ws.Range("I22").Value = 553.2
$ws.Range("K22").Value = 553.2
$ws.Range("M22").Value = -380.2

$ws.Range("H80").Value = 4004.5625
$ws.Range("J80").Value = 1181
$ws.Range("L80").Value = 1181
$ws.Range("N80").Value = -3177

$ws.Range("H83").Value = 4004.5625
$ws.Range("J83").Value = 1181
$ws.Range("L83").Value = 5905
$ws.Range("N83").Value = -15889

$ws.Range("H86").Value = 2845.0908
$ws.Range("I86").Value = 2831.7778
$ws.Range("J86").Value = 2905
$ws.Range("K86").Value = 2831.7778
$ws.Range("L86").Value = 2905
$ws.Range("M86").Value = -1708.7778
$ws.Range("N86").Value = -5151

$ws.Range("H89").Value = 2845.0908
$ws.Range("I89").Value = 2831.7778
$ws.Range("J89").Value = 2905
$ws.Range("K89").Value = 14158.889
$ws.Range("L89").Value = 14525
$ws.Range("M89").Value = -8542.888999999999
$ws.Range("N89").Value = -25757

$ws.Range("H134").Value = 6901.2173
$ws.Range("I134").Value = 6901.2173
$ws.Range("K134").Value = 20703.6519
$ws.Range("M134").Value = -18168.6519

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 37285.715
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250

$ws.Range("H59").Value = 64500
$ws.Range("J59").Value = 75000
$ws.Range("L59").Value = 75000
$ws.Range("N59").Value = -77290

$ws.Range("H60").Value = 17621.143
$ws.Range("J60").Value = 15000
$ws.Range("L60").Value = 15000
$ws.Range("N60").Value = -16022

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H122").Value = 4461.5557
$ws.Range("I122").Value = 4033.4
$ws.Range("J122").Value = 4996.75
$ws.Range("K122").Value = 12100.2
$ws.Range("L122").Value = 14990.25
$ws.Range("M122").Value = -9650.200000000001
$ws.Range("N122").Value = -19890.25

$ws.Range("H134").Value = 55001.5
$ws.Range("I134").Value = 55001.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 165004.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -162469.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 122.4
$ws.Range("I2").Value = 64.8
$ws.Range("J2").Value = 180
$ws.Range("K2").Value = 388.8
$ws.Range("L2").Value = 1080
$ws.Range("M2").Value = -275.8
$ws.Range("N2").Value = -1306

$ws.Range("H12").Value = 129.86667
$ws.Range("I12").Value = 72.75
$ws.Range("K12").Value = 218.25
$ws.Range("M12").Value = -45.25

$ws.Range("H36").Value = 858
$ws.Range("I36").Value = 260.66666
$ws.Range("J36").Value = 2650
$ws.Range("K36").Value = 781.9999799999999
$ws.Range("L36").Value = 7950
$ws.Range("M36").Value = -612.9999799999999
$ws.Range("N36").Value = -8288

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4118.8
$ws.Range("I122").Value = 3671
$ws.Range("J122").Value = 10388
$ws.Range("K122").Value = 11013
$ws.Range("L122").Value = 31164
$ws.Range("M122").Value = -8563
$ws.Range("N122").Value = -36064

$ws.Range("H132").Value = 4399.3335
$ws.Range("I132").Value = 2849
$ws.Range("K132").Value = 8547
$ws.Range("M132").Value = -6017

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 920.7273
$ws.Range("I22").Value = 914.8889
$ws.Range("K22").Value = 914.8889
$ws.Range("M22").Value = -619.8889

$ws.Range("H27").Value = 920.7273
$ws.Range("I27").Value = 914.8889
$ws.Range("K27").Value = 914.8889
$ws.Range("M27").Value = -807.8889

$ws.Range("H39").Value = 22495
$ws.Range("J39").Value = 22495
$ws.Range("L39").Value = 22495
$ws.Range("N39").Value = -23415

$ws.Range("H46").Value = 3384.65
$ws.Range("I46").Value = 2850
$ws.Range("J46").Value = 3919.3
$ws.Range("K46").Value = 2850
$ws.Range("L46").Value = 3919.3
$ws.Range("M46").Value = -2662
$ws.Range("N46").Value = -4295.3

$ws.Range("H61").Value = 7581.353
$ws.Range("J61").Value = 1633
$ws.Range("L61").Value = 1633
$ws.Range("N61").Value = -2037

$ws.Range("H68").Value = 13337686
$ws.Range("I68").Value = 14495093
$ws.Range("J68").Value = 27500
$ws.Range("K68").Value = 14495093
$ws.Range("L68").Value = 27500
$ws.Range("M68").Value = -14494344
$ws.Range("N68").Value = -28998

$ws.Range("H71").Value = 13337686
$ws.Range("I71").Value = 14495093
$ws.Range("J71").Value = 27500
$ws.Range("K71").Value = 72475465
$ws.Range("L71").Value = 137500
$ws.Range("M71").Value = -72471721
$ws.Range("N71").Value = -144988

$ws.Range("H113").Value = 7581.353
$ws.Range("J113").Value = 1633
$ws.Range("L113").Value = 1633
$ws.Range("N113").Value = -5973

$ws.Range("H132").Value = 30471.541
$ws.Range("I132").Value = 36669.633
$ws.Range("K132").Value = 110008.899
$ws.Range("M132").Value = -107478.899

$ws.Range("H136").Value = 5822827
$ws.Range("I136").Value = 8588340
$ws.Range("J136").Value = 15248.5
$ws.Range("K136").Value = 25765020
$ws.Range("L136").Value = 45745.5
$ws.Range("M136").Value = -25762470
$ws.Range("N136").Value = -50845.5

$ws.Range("H138").Value = 127606.5
$ws.Range("J138").Value = 127606.5
$ws.Range("L138").Value = 127606.5
$ws.Range("N138").Value = -137886.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 10749.5
$ws.Range("J52").Value = 12666
$ws.Range("L52").Value = 12666
$ws.Range("N52").Value = -13118

$ws.Range("H62").Value = 166677090
$ws.Range("I62").Value = 166677090
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 166677090
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -166676466
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 166677090
$ws.Range("I65").Value = 166677090
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 166677090
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -833382330
$ws.Range("N65").ClearContents()

$ws.Range("H104").Value = 21541
$ws.Range("J104").Value = 21541
$ws.Range("L104").Value = 21541
$ws.Range("N104").Value = -28529

$ws.Range("H113").Value = 628.1177
$ws.Range("I113").Value = 376.3846
$ws.Range("K113").Value = 1129.1538
$ws.Range("M113").Value = 1040.8462

$ws.Range("H122").Value = 9103.870999999999
$ws.Range("I122").Value = 6835.8096
$ws.Range("J122").Value = 13866.8
$ws.Range("K122").Value = 20507.4288
$ws.Range("L122").Value = 41600.39999999999
$ws.Range("M122").Value = -18057.4288
$ws.Range("N122").Value = -46500.39999999999

$ws.Range("H126").Value = 8743.9375
$ws.Range("I126").Value = 4541.727
$ws.Range("K126").Value = 13625.181
$ws.Range("M126").Value = -11155.181
